$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H76").Value = 108954.95
$ws.Range("I76").Value = 150944.06
$ws.Range("J76").Value = 3982.1667
$ws.Range("K76").Value = 150944.06
$ws.Range("L76").Value = 3982.1667
$ws.Range("M76").Value = -150629.06
$ws.Range("N76").Value = -4612.1667
$ws.Range("H79").Value = 108954.95
$ws.Range("I79").Value = 150944.06
$ws.Range("J79").Value = 3982.1667
$ws.Range("K79").Value = 150944.06
$ws.Range("L79").Value = 3982.1667
$ws.Range("M79").Value = -149852.06
$ws.Range("N79").Value = -6166.1667
$ws.Range("H98").Value = 760.69696
$ws.Range("I98").Value = 451.67856
$ws.Range("J98").Value = 2491.2
$ws.Range("K98").Value = 451.67856
$ws.Range("L98").Value = 2491.2
$ws.Range("M98").Value = 1046.32144
$ws.Range("N98").Value = -5487.2
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 760.69696
$ws.Range("I122").Value = 451.67856
$ws.Range("J122").Value = 2491.2
$ws.Range("K122").Value = 1355.03568
$ws.Range("L122").Value = 7473.599999999999
$ws.Range("M122").Value = 1094.96432
$ws.Range("N122").Value = -12373.6
$ws.Range("H129").Value = 7253.5
$ws.Range("I129").Value = 452.22223
$ws.Range("K129").Value = 1356.66669
$ws.Range("M129").Value = 3643.33331
$ws.Range("H138").Value = 2861.2056
$ws.Range("I138").Value = 1712.381
$ws.Range("K138").Value = 5137.143
$ws.Range("M138").Value = 2.856999999999971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 2126.5
$ws.Range("I16").Value = 2126.5
$ws.Range("K16").Value = 2126.5
$ws.Range("M16").Value = -1839.5
$ws.Range("H26").Value = 2503.5
$ws.Range("I26").Value = 2503.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 2503.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -2173.5
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5043.804
$ws.Range("I94").Value = 583.84375
$ws.Range("J94").Value = 15238
$ws.Range("K94").Value = 583.84375
$ws.Range("L94").Value = 15238
$ws.Range("M94").Value = -132.84375
$ws.Range("N94").Value = -16140
$ws.Range("H99").Value = 1704.579
$ws.Range("I99").Value = 1490.6666
$ws.Range("J99").Value = 2071.2856
$ws.Range("K99").Value = 1490.6666
$ws.Range("L99").Value = 2071.2856
$ws.Range("M99").Value = 7.333399999999983
$ws.Range("N99").Value = -5067.2856
$ws.Range("H105").Value = 2362.9412
$ws.Range("I105").Value = 1825.9
$ws.Range("J105").Value = 3130.1428
$ws.Range("K105").Value = 1825.9
$ws.Range("L105").Value = 3130.1428
$ws.Range("M105").Value = -78.90000000000009
$ws.Range("N105").Value = -6624.1428
$ws.Range("H134").Value = 1467.3043
$ws.Range("I134").Value = 973.41174
$ws.Range("J134").Value = 2866.6667
$ws.Range("K134").Value = 2920.23522
$ws.Range("L134").Value = 8600.000100000001
$ws.Range("M134").Value = -385.23522
$ws.Range("N134").Value = -13670.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 622.5
$ws.Range("I8").Value = 600
$ws.Range("J8").Value = 630
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 630
$ws.Range("M8").Value = -460
$ws.Range("N8").Value = -910
$ws.Range("H58").Value = 1457.8667
$ws.Range("I58").Value = 811.4286
$ws.Range("J58").Value = 2023.5
$ws.Range("K58").Value = 811.4286
$ws.Range("L58").Value = 2023.5
$ws.Range("M58").Value = -608.4286
$ws.Range("N58").Value = -2429.5
$ws.Range("H136").Value = 1457.8667
$ws.Range("I136").Value = 811.4286
$ws.Range("J136").Value = 2023.5
$ws.Range("K136").Value = 2434.2858
$ws.Range("L136").Value = 6070.5
$ws.Range("M136").Value = 115.7142000000003
$ws.Range("N136").Value = -11170.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 485256.22
$ws.Range("I5").Value = 531.6087
$ws.Range("J5").Value = 1042689.56
$ws.Range("K5").Value = 1594.8261
$ws.Range("L5").Value = 3128068.68
$ws.Range("M5").Value = -1482.8261
$ws.Range("N5").Value = -3128292.68
$ws.Range("H22").Value = 3400.1428
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 3070.1
$ws.Range("K22").Value = 30003
$ws.Range("L22").Value = 9210.299999999999
$ws.Range("M22").Value = -29834
$ws.Range("N22").Value = -9548.299999999999
$ws.Range("H27").Value = 3400.1428
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 3070.1
$ws.Range("K27").Value = 30003
$ws.Range("L27").Value = 9210.299999999999
$ws.Range("M27").Value = -29901
$ws.Range("N27").Value = -9414.299999999999
$ws.Range("H107").Value = 338.77777
$ws.Range("I107").Value = 247
$ws.Range("J107").Value = 522.3333
$ws.Range("K107").Value = 741
$ws.Range("L107").Value = 1566.9999
$ws.Range("M107").Value = 1179
$ws.Range("N107").Value = -5406.9999
$ws.Range("H113").Value = 614.25
$ws.Range("I113").Value = 394.81818
$ws.Range("J113").Value = 847.8387
$ws.Range("K113").Value = 1184.45454
$ws.Range("L113").Value = 2543.5161
$ws.Range("M113").Value = 985.54546
$ws.Range("N113").Value = -6883.5161
$ws.Range("H120").Value = 17290
$ws.Range("I120").Value = 18720
$ws.Range("J120").Value = 13000
$ws.Range("K120").Value = 56160
$ws.Range("L120").Value = 39000
$ws.Range("M120").Value = -51322
$ws.Range("N120").Value = -48676
$ws.Range("H122").Value = 720.9677
$ws.Range("I122").Value = 297.21054
$ws.Range("J122").Value = 1391.9166
$ws.Range("K122").Value = 2674.89486
$ws.Range("L122").Value = 12527.2494
$ws.Range("M122").Value = -224.8948599999999
$ws.Range("N122").Value = -17427.2494
$ws.Range("H132").Value = 990.39026
$ws.Range("I132").Value = 715.4286
$ws.Range("J132").Value = 1047
$ws.Range("K132").Value = 6438.8574
$ws.Range("L132").Value = 9423
$ws.Range("M132").Value = -3908.8574
$ws.Range("N132").Value = -14483
$ws.Range("H135").Value = 485256.22
$ws.Range("I135").Value = 531.6087
$ws.Range("J135").Value = 1042689.56
$ws.Range("K135").Value = 4784.4783
$ws.Range("L135").Value = 9384206.040000001
$ws.Range("M135").Value = -2249.4783
$ws.Range("N135").Value = -9389276.040000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 809.5
$ws.Range("I16").Value = 842.2857
$ws.Range("J16").Value = 694.75
$ws.Range("K16").Value = 842.2857
$ws.Range("L16").Value = 694.75
$ws.Range("M16").Value = -672.2857
$ws.Range("N16").Value = -1034.75
$ws.Range("H105").Value = 42980
$ws.Range("J105").Value = 42980
$ws.Range("L105").Value = 42980
$ws.Range("N105").Value = -49968
$ws.Range("H136").Value = 19616404
$ws.Range("I136").Value = 66691350
$ws.Range("J136").Value = 1841.6666
$ws.Range("K136").Value = 200074050
$ws.Range("L136").Value = 5524.9998
$ws.Range("M136").Value = -200071500
$ws.Range("N136").Value = -10624.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2000000
$ws.Range("I5").Value = 2000000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2000000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1999888
$ws.Range("N5").ClearContents()
$ws.Range("H119").Value = 191883.33
$ws.Range("J119").Value = 191883.33
$ws.Range("L119").Value = 191883.33
$ws.Range("N119").Value = -201559.33
$ws.Range("H136").Value = 6967163.5
$ws.Range("I136").Value = 13236001
$ws.Range("J136").Value = 1787.7778
$ws.Range("K136").Value = 13236001
$ws.Range("L136").Value = 5363.3334
$ws.Range("M136").Value = -39705453
$ws.Range("N136").Value = -10463.3334
